$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.278.55'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.58%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.423.88'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.06%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '563.47'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.14%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.64%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.422.22'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.94%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.79%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.35%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.74%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.76%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +5.14%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.861.68'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.075.64'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.19%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.425.69'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.94%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +3.20%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.50'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.91%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.73%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.62'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.99%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.93'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.96%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '583.95'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +12.86%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.542.97'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0944'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.90%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.46'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +5.02%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.26'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.85%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.49%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.88'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.55'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.38%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.72'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.76%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.79'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.30%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.58%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '153.17'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.67'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.96%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.36%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.33'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +7.77%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '150.30'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.69%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.51%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.35'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.69%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.16%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.96%  '
